$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    "DIVERSOS",
    "EQUIPAMENTOS E ARTIGOS PARA USO MÉDICO, DENTÁRIO E VETERINÁRIO",
    "INSTRUMENTOS E EQUIPAMENTOS DE LABORATÓRIO",
    "MATERIAIS MANUFATURADOS, NÃO METÁLICOS",
    "SUBSISTÊNCIA",
    "SUBSTÂNCIAS E PRODUTOS QUÍMICOS"
)

$template = $ws.Cells.Item(2, 1)

$row = 3
foreach ($g in $groups) {
    $cell = $ws.Cells.Item($row, 1)
    $template.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $g
    $cell.RowHeight = 18.75
    $row++
}
